# Applies the "rolling future window + loss over future" edit:
#  - Slide 1: update meeting date
#  - Slide 4: update several numbers in the results table + highlight some cells
#  - Slide 5: shorten the "Autoregressive models..." bullet text

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 - Title slide: bump the meeting date
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(2)
$tr = $subtitle.TextFrame.TextRange
$fullText = $tr.Text
$oldDate = "Meeting: 20/06/2023"
$newDate = "Meeting: 04/07/2023"
$idx = $fullText.IndexOf($oldDate)
if ($idx -ge 0) {
    $run = $tr.Characters($idx + 1, $oldDate.Length)
    $run.Text = $newDate
}

# ---------------------------------------------------------------------------
# Slide 4 - Results table: refresh the "All Action" numbers (rolling future
# window results), and highlight the "Teacher Forcing Feeding Split" figures
# that correspond to the new loss-over-future comparison.
# ---------------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$tableShape = $slide4.Shapes.Item(2)
$tbl = $tableShape.Table

function Set-CellText($row, $col, $text) {
    $cell = $tbl.Cell($row, $col)
    $cell.Shape.TextFrame.TextRange.Text = $text
}

function Add-CellHighlight($row, $col) {
    $cell = $tbl.Cell($row, $col)
    $cell.Shape.TextFrame.TextRange.Font.Highlight.RGB = 65535
}

# "All Action" column (6) value refresh
Set-CellText 2 6 "0.4267"
Set-CellText 4 6 "0.4303"
Set-CellText 5 6 "0.0857"
Set-CellText 6 6 "0.4329"
Set-CellText 7 6 " 0.6908"
Set-CellText 8 6 "0.4349"
Set-CellText 9 6 "0.0895"

# "Baseline" column (3) value refresh
Set-CellText 9 3 "0.1501"

# Highlight the "Teacher Forcing Feeding Split" column (5) on rows that now
# mark the best score for that metric
Add-CellHighlight 3 5
Add-CellHighlight 6 5
Add-CellHighlight 8 5
Add-CellHighlight 9 5

# ---------------------------------------------------------------------------
# Slide 5 - Future slide: shorten the autoregressive-models bullet
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$body = $slide5.Shapes.Item(2)
$btr = $body.TextFrame.TextRange
$bFull = $btr.Text
$oldBullet = "Autoregressive models but not rolling window:"
$newBullet = "Autoregressive models:"
$bidx = $bFull.IndexOf($oldBullet)
if ($bidx -ge 0) {
    $brun = $btr.Characters($bidx + 1, $oldBullet.Length)
    $brun.Text = $newBullet
}
